# "updated user management page"
# Adds a new "UserDetails" worksheet (after the existing "logins" sheet)
# that holds a simple key/value table of a user's profile details, and
# makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Start from a copy of the existing "logins" sheet so the new sheet
# inherits the same sheet-level formatting (outline levels, etc.) as the
# rest of the workbook, then re-purpose it as "UserDetails".
$logins = $wb.Worksheets.Item("logins")
$logins.Copy($null, $logins)
$ws = $wb.ActiveSheet
$ws.Name = "UserDetails"

# Clear any inherited data from the copy before writing the new content.
$ws.Cells.Clear()

$rows = @(
    @("Username", "Andreson"),
    @("Password", "anil123"),
    @("EmployeeName", "Kevin  Mathews"),
    @("Status", "Enabled"),
    @("UserRole", "ESS")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Match the column sizing / page setup of the shipped sheet.
$ws.Columns.Item(1).ColumnWidth = 26.65

$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Select B1 on the new sheet, which is left as the active sheet/tab.
$ws.Range("B1").Select() | Out-Null
